$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Differentiate between 3D and 1D: update the gradient value in D2,
# make D3 reference D2 via a formula, and let dependent formulas recalc.
$ws.Range("D2").Value = 1
$ws.Range("D3").Formula = "=D2"

# Update the active cell selection as recorded in the saved file.
$ws.Range("D4").Select()

$wb.Save()
